$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose target text would otherwise be auto-converted to a number by Excel;
# force them to remain as text (matching the original inlineStr/text cells),
# then clear the temporary number-format override so no stray style sticks around.
$textForceCells = @("D5", "D6", "D14", "D20", "D21", "D22", "D23", "D25", "D26", "D27", "D32", "D37", "D38", "D39", "D41", "D46", "D47")
foreach ($addr in $textForceCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D5").Value = '574.77'
$ws.Range("D6").Value = '147.04'
$ws.Range("D14").Value = '29.01'
$ws.Range("D20").Value = '11.02'
$ws.Range("D21").Value = '326.47'
$ws.Range("D22").Value = '2.22'
$ws.Range("D23").Value = '4.13'
$ws.Range("D25").Value = '10.08'
$ws.Range("D26").Value = '65.57'
$ws.Range("D27").Value = '645.00'
$ws.Range("D32").Value = '7.99'
$ws.Range("D37").Value = '4.74'
$ws.Range("D38").Value = '2.82'
$ws.Range("D39").Value = '152.15'
$ws.Range("D41").Value = '18.71'
$ws.Range("D46").Value = '152.59'
$ws.Range("D47").Value = '15.28'

foreach ($addr in $textForceCells) {
    $ws.Range($addr).ClearFormats()
}

# Remaining cells: plain text assignment is sufficient (Excel keeps them as text).
$ws.Range("D2").Value = '62.817.89'
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").Value = '2.463.09'
$ws.Range("E3").Value = '  +0.72%  '
$ws.Range("E4").Value = '  +0.06%  '
$ws.Range("E5").Value = '  -0.22%  '
$ws.Range("E6").Value = '  +0.97%  '
$ws.Range("E7").Value = '  -0.04%  '
$ws.Range("E8").Value = '  -0.91%  '
$ws.Range("D9").Value = '2.462.58'
$ws.Range("E10").Value = '  +0.77%  '
$ws.Range("E11").Value = '  -0.23%  '
$ws.Range("E12").Value = '  +0.15%  '
$ws.Range("E13").Value = '  +0.99%  '
$ws.Range("E14").Value = '  +2.23%  '
$ws.Range("E15").Value = '  -0.41%  '
$ws.Range("D16").Value = '2.910.20'
$ws.Range("E16").Value = '  +0.69%  '
$ws.Range("D17").Value = '62.780.46'
$ws.Range("E17").Value = '  +0.20%  '
$ws.Range("D18").Value = '2.465.27'
$ws.Range("E18").Value = '  +0.89%  '
$ws.Range("E19").Value = '  +0.21%  '
$ws.Range("E20").Value = '  +0.14%  '
$ws.Range("E21").Value = '  -1.11%  '
$ws.Range("B22").Value = 'SuiNetwork'
$ws.Range("C22").Value = 'https://coinranking.com/coin/3xJluUMvp+suinetwork-sui'
$ws.Range("E22").Value = '  +8.86%  '
$ws.Range("B23").Value = 'Polkadot'
$ws.Range("C23").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("E23").Value = '  -0.07%  '
$ws.Range("E24").Value = '  -0.09%  '
$ws.Range("E25").Value = '  +17.99%  '
$ws.Range("E26").Value = '  -1.09%  '
$ws.Range("E27").Value = '  -0.30%  '
$ws.Range("D28").Value = '0.0₃0985'
$ws.Range("E28").Value = '  -0.40%  '
$ws.Range("E30").Value = '  -15.32%  '
$ws.Range("E31").Value = '  -0.47%  '
$ws.Range("E32").Value = '  -2.52%  '
$ws.Range("E33").Value = '  -1.24%  '
$ws.Range("E34").Value = '  -2.80%  '
$ws.Range("E36").Value = '  +2.82%  '
$ws.Range("E37").Value = '  -0.24%  '
$ws.Range("E38").Value = '  +3.81%  '
$ws.Range("E39").Value = '  -0.81%  '
$ws.Range("E40").Value = '  -1.38%  '
$ws.Range("E41").Value = '  -0.17%  '
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("E44").Value = '  -36.92%  '
$ws.Range("E45").Value = '  +0.00%  '
$ws.Range("E46").Value = '  +5.10%  '
$ws.Range("E47").Value = '  +2.23%  '
$ws.Range("E48").Value = '  -1.30%  '
$ws.Range("E49").Value = '  -0.51%  '
$ws.Range("E50").Value = '  +0.34%  '
$ws.Range("E51").Value = '  -1.03%  '
